$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header label in B2 (was "unnamed: 1_level_1", should read "total")
$ws.Range("B2").Value = "total"

# Remove the now-unneeded section header rows entirely (delete from the bottom up
# so row indices for the earlier deletion stay valid):
#   row 8 -> "grandes regiões e unidades da federação" (empty section header row)
#   row 5 -> "situação do domicílio" (empty section header row)
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
